$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 174 - this shifts rows 174:212 down to 175:213
# and preserves styles/formatting of the shifted cells (including the date
# style on column D).
$ws.Rows.Item(174).Insert()

# Populate the newly inserted row 174 with the new record's data.
$ws.Cells.Item(174, 1).Value = 6
$ws.Cells.Item(174, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(174, 3).Value = "Metropolitana"
$ws.Cells.Item(174, 4).Value = 44722
$ws.Cells.Item(174, 5).Value = 13
$ws.Cells.Item(174, 6).Value = 100112029
$ws.Cells.Item(174, 7).Value = "Orégano"
$ws.Cells.Item(174, 8).Value = "Sin especificar"
$ws.Cells.Item(174, 9).Value = "Primera"
$ws.Cells.Item(174, 10).Value = 47
$ws.Cells.Item(174, 11).Value = 12000
$ws.Cells.Item(174, 12).Value = 13000
$ws.Cells.Item(174, 13).Value = 12468
$ws.Cells.Item(174, 14).Value = "$/docena de atados"
$ws.Cells.Item(174, 15).Value = "Región Metropolitana"
$ws.Cells.Item(174, 16).Value = 4156
$ws.Cells.Item(174, 17).Value = 3
$ws.Cells.Item(174, 18).Value = "Hortaliza"
